# Added two new Mac-Addresses: append two rows (31 & 32) to the
# master-reg_center_user_machine_h data sheet, following the exact
# pattern of the existing rows (e.g. row 30).
#
# Columns: A=regcntr_id  B=usr_id  C=machine_id  D=lang_code
#          E=is_active   F=cr_by   G=cr_dtimes    H=eff_dtimes

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @{ Row = 31; RegCntrId = 10001; UsrId = 110030; MachineId = 10030 },
    @{ Row = 32; RegCntrId = 10001; UsrId = 110031; MachineId = 10031 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.RegCntrId
    $ws.Cells.Item($row, 2).Value = $r.UsrId
    $ws.Cells.Item($row, 3).Value = $r.MachineId
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
    $ws.Cells.Item($row, 8).Value = "now()"
}

# Scroll the view down to the newly added rows and leave the last new
# row's cr_by cell selected, matching the author's final cursor position.
$ws.Range("A25").Select() | Out-Null
$ws.Range("F30").Select() | Out-Null
